$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row for LeetCode problem #1 "Two Sum(s)" -------------------------------
$newRow = 19

# Copy the formatting (fonts, borders, alignment, fill, number format) of the
# row above (row 18, the last existing data row) into the new row so the new
# entry visually matches the rest of the table.
$ws.Range("A18:I18").Copy() | Out-Null
$ws.Range("A19:I19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Rows.Item($newRow).RowHeight = 75

$ws.Range("A19").Value = 1
$ws.Range("B19").Value = "Two Sums"
$ws.Range("C19").Value = "Easy"
$ws.Range("D19").Value = "Arrays/HashMap"
$ws.Range("E19").Value = "Accepted"
$ws.Range("F19").Value = "O(n)"
$ws.Range("G19").Value = "O(n)"
$ws.Range("H19").Value = "Phase 3-4"
$ws.Range("I19").Value = "***This is an important interview question as well. Because it is between x=target -y; inserting the elements into HashMap and compare from there. It seems that second part is two pointers but need to figure out how to do binary search first, since second part will always be sorted."

# Keep the view pointed at the newly edited area, matching where the author
# left the cursor after adding the row.
$ws.Application.Goto($ws.Range("A7"), $true)
$ws.Range("H12").Select() | Out-Null
